$wb = $excel.ActiveWorkbook

# Replace the raw CDISC-coded values with separated, cleaner representations
# (code kept on its own, label/term kept on its own) for studyType/studyPhase.
$ws = $wb.Worksheets.Item("study")
$ws.Range("C2").Value = "Interventional Study"
$ws.Range("D2").Value = "C15602"

# Make "study" the active sheet/tab with the same selection state captured
# in the saved file (cursor left on C11).
$ws.Activate()
$ws.Range("C11").Select() | Out-Null
